$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: ATRIBUTO "tipo de dato" tag, "NUM" for every attribute row ---
# --- (2-18, the numeric attributes).                                     ---
foreach ($r in 2..18) {
    $ws.Cells.Item($r, 3).Value = "NUM"
}

# --- Column D: validation notes (new shared strings, added in diff order ---
# --- before the row-19 "cateogorica int" tag).                           ---
$ws.Cells.Item(4, 4).Value = "Revisar que este en el rango de 0 a 1"

# Row 9 note spans (and is centered/wrapped across) the merged block D9:D12.
$note = $ws.Cells.Item(9, 4)
$note.Value = "Revisar que este en el rango de 0 a 1"
$note.HorizontalAlignment = -4108
$note.VerticalAlignment = -4108
$note.WrapText = $true

# Copy that same formatting down onto D10:D12 so every cell of the merge
# block shares one style entry, then merge the block.
$note.Copy()
$ws.Range("D10:D12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D9:D12").Merge()

# --- Last attribute (MESES_CLIENTE) is flagged as a categorical integer. ---
$ws.Cells.Item(19, 3).Value = "cateogorica int"

# --- Column widths / row heights refreshed alongside the new column.     ---
$ws.Columns.Item(4).ColumnWidth = 14.498697916666666
foreach ($r in 1..19) {
    $ws.Rows.Item($r).RowHeight = 25
}

# --- Selection cursor moved to A18 on last save. ---
$ws.Range("A18").Select() | Out-Null
